$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.0345785
$ws.Range("H2").Value = 0.069157
$ws.Range("M2").Value = 18.296032
$ws.Range("N2").Value = 54.888096
$ws.Range("O2").Value = 0.1408813274592244
$ws.Range("P2").Value = 0.1447879988537689
$ws.Range("Q2").Value = 0.632649342512
$ws.Range("R2").Value = 3.795896055072
$ws.Range("S2").Value = 0.1408813274592244
$ws.Range("T2").Value = 0.1447879988537689

# Row 3
$ws.Range("G3").Value = 0.0345785
$ws.Range("H3").Value = 0.069157
$ws.Range("O3").Value = 0.2147385011396954
$ws.Range("P3").Value = 0.2206932488329456
$ws.Range("Q3").Value = 0.9643163789563332
$ws.Range("R3").Value = 5.785898273737999
$ws.Range("S3").Value = 0.2147385011396954
$ws.Range("T3").Value = 0.2206932488329456

# Row 4
$ws.Range("G4").Value = 0.0345785
$ws.Range("H4").Value = 0.069157
$ws.Range("M4").Value = 32.05318933333334
$ws.Range("N4").Value = 96.15956800000001
$ws.Range("O4").Value = 0.2468128533324521
$ws.Range("P4").Value = 0.2536570301393388
$ws.Range("Q4").Value = 1.108351207362667
$ws.Range("R4").Value = 6.650107244176001
$ws.Range("S4").Value = 0.2468128533324521
$ws.Range("T4").Value = 0.2536570301393388

# Row 5
$ws.Range("G5").Value = 0.0345785
$ws.Range("H5").Value = 0.069157
$ws.Range("M5").Value = 10.5123315
$ws.Range("N5").Value = 21.024663
$ws.Range("O5").Value = 0.08094603334818277
$ws.Range("P5").Value = 0.05546045689660794
$ws.Range("Q5").Value = 0.36350065477275
$ws.Range("R5").Value = 1.454002619091
$ws.Range("S5").Value = 0.08094603334818277
$ws.Range("T5").Value = 0.05546045689660794

# Row 6
$ws.Range("G6").Value = 0.0345785
$ws.Range("H6").Value = 0.069157
$ws.Range("M6").Value = 41.11909833333333
$ws.Range("N6").Value = 123.357295
$ws.Range("O6").Value = 0.3166212847204454
$ws.Range("P6").Value = 0.3254012652773388
$ws.Range("Q6").Value = 1.421836741719167
$ws.Range("R6").Value = 8.531020450314999
$ws.Range("S6").Value = 0.3166212847204454
$ws.Range("T6").Value = 0.3254012652773388
